# Update the mantel correlogram table values after recreating the
# correlograms with Euclidean distances.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Columns: 1 = Distance Class (m), 2 = N, 3 = Mantel r, 4 = p
# Rows: 1 = header, 2 = "5,000", 3 = "15,000", 4 = "25,000"

function Set-CellValue($row, $col, $oldValue, $newValue) {
    $cell = $table.Cell($row, $col)
    $current = $cell.Range.Text
    if ($current -notmatch [regex]::Escape($oldValue)) {
        throw "Unexpected content in row $row col $col : expected '$oldValue' but found '$current'"
    }
    $cell.Range.Text = $newValue
    return $cell
}

# Row "5,000": Mantel r and p
Set-CellValue 2 3 "-0.047" "0.029"
Set-CellValue 2 4 "0.093" "0.192"

# Row "15,000": Mantel r and p
Set-CellValue 3 3 "-0.002" "0.035"
Set-CellValue 3 4 "0.46" "0.252"

# Row "25,000": Mantel r and p
Set-CellValue 4 3 "0.073" "-0.007"
$pCell = Set-CellValue 4 4 "0.042" "0.45"

# The last cell (p-value for the 25,000 row) was bolded before; it is no
# longer significant/highlighted, so remove the bold formatting.
$pCell.Range.Font.Bold = $false
